$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.98"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.202"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05742"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.477"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.241"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8152"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8688"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1371"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06934"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03166"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02872"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09327"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.813"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001527"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.01018"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006148"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001235"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00008696"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.582"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.156"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1330"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002328"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006227"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1050"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001983"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008386"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005442"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4540"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003323"
